# Update cryptocurrency price (D) and hourly volume change (E) columns
# for rows 2-51 on Sheet1, per the Tue Jul 25 10:13:54 UTC 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Volume(1h)) --------------------------------------------
# These are plain percentage strings (never numeric-looking), so a direct
# .Value assignment keeps them as text with their original formatting.
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  +4.20%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E12").Value = '  +5.40%  '
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("E16").Value = '  -2.10%  '
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  -3.27%  '
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("E30").Value = '  +0.98%  '
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("E41").Value = '  +6.00%  '
$ws.Range("E42").Value = '  +6.66%  '
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("E51").Value = '  -3.84%  '

# --- Column D (Price) ---------------------------------------------------
# Many of the new price strings parse as plain numbers (e.g. "1.001",
# "0.07700"), which a direct .Value assignment would silently convert to
# a Double and mangle (dropping trailing zeros / turning "1.000" into 1).
# To preserve the exact text, stage each value as a literal-text formula
# in a scratch column, then Copy / PasteSpecial(xlPasteValues) it onto the
# Price column, which keeps the cell a plain text value. The scratch
# column is cleared afterwards.
$ws.Range("G2").Formula = '="29.195.01"'
$ws.Range("G3").Formula = '="1.857.19"'
$ws.Range("G4").Formula = '="1.000"'
$ws.Range("G5").Formula = '="238.16"'
$ws.Range("G6").Formula = '="0.6894"'
$ws.Range("G7").Formula = '="1.001"'
$ws.Range("G8").Formula = '="0.07700"'
$ws.Range("G9").Formula = '="0.3052"'
$ws.Range("G10").Formula = '="23.26"'
$ws.Range("G11").Formula = '="0.08061"'
$ws.Range("G12").Formula = '="1.979.57"'
$ws.Range("G13").Formula = '="0.7218"'
$ws.Range("G14").Formula = '="5.193"'
$ws.Range("G15").Formula = '="89.36"'
$ws.Range("G16").Formula = '="29.198.28"'
$ws.Range("G17").Formula = '="5.734"'
$ws.Range("G18").Formula = '="0.000007790"'
$ws.Range("G19").Formula = '="13.25"'
$ws.Range("G20").Formula = '="235.02"'
$ws.Range("G21").Formula = '="1.000"'
$ws.Range("G22").Formula = '="2.108.63"'
$ws.Range("G23").Formula = '="1.001"'
$ws.Range("G24").Formula = '="7.453"'
$ws.Range("G25").Formula = '="161.90"'
$ws.Range("G26").Formula = '="8.958"'
$ws.Range("G27").Formula = '="0.1430"'
$ws.Range("G28").Formula = '="18.06"'
$ws.Range("G29").Formula = '="1.953"'
$ws.Range("G30").Formula = '="1.400"'
$ws.Range("G31").Formula = '="4.534"'
$ws.Range("G32").Formula = '="1.486"'
$ws.Range("G33").Formula = '="4.012"'
$ws.Range("G34").Formula = '="0.05194"'
$ws.Range("G35").Formula = '="1.184"'
$ws.Range("G36").Formula = '="0.7046"'
$ws.Range("G37").Formula = '="1.022"'
$ws.Range("G38").Formula = '="2.674"'
$ws.Range("G39").Formula = '="0.01849"'
$ws.Range("G40").Formula = '="2.678"'
$ws.Range("G41").Formula = '="0.9287"'
$ws.Range("G42").Formula = '="1.095.93"'
$ws.Range("G43").Formula = '="5.973"'
$ws.Range("G44").Formula = '="0.4287"'
$ws.Range("G45").Formula = '="71.00"'
$ws.Range("G46").Formula = '="1.001"'
$ws.Range("G47").Formula = '="102.57"'
$ws.Range("G48").Formula = '="1.794"'
$ws.Range("G49").Formula = '="2.004.92"'
$ws.Range("G50").Formula = '="9.152"'
$ws.Range("G51").Formula = '="7.001"'

$scratch = $ws.Range("G2:G51")
$scratch.Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

Write-Output "Updated D2:D51 and E2:E51"
